$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "T1234567J, T1234567J"
$ws.Range("N2").Value = "T1234567J, T1234567J, T1234567J"
$ws.Range("N2").Value = "T1234567J, T1234567J, T1234567J, T1234567J"
